# Daily attendance processing - 2025-11-21 10:49:03
#
# The "Recorded By" column (G) stores a comma-separated list of the users
# who touched the attendance record (e.g. "dnasr281@gmail.com, System").
# Going forward "System" should be listed first in that list, so every
# existing row whose list has more than one entry, includes the literal
# token "System", and does not already have "System" in the first
# position gets its entries put in reverse order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G - "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $parts = $value -split ", "
    $count = $parts.Length

    if ($count -gt 1) {
        $firstPart = $parts[0]
        $lastPart = $parts[$count - 1]

        # Only reorder when the (case-sensitive) token "System" is the last
        # entry and not already the first entry.
        if ($lastPart.Equals("System") -and -not $firstPart.Equals("System")) {
            $reversedParts = $parts[($count - 1)..0]
            $cell.Value2 = $reversedParts -join ", "
        }
    }
}
